$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1109
$wsExhibit.Range("F4").Value = 1778
$wsExhibit.Range("F6").Value = 285

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1109
$wsAll.Range("F4").Value = 1778
$wsAll.Range("F7").Value = 285
